$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading numbers: merge split run "3." + "<digit>" into a single "3.<digit>"
#    run (no visible text or formatting change - Word just coalesces the runs
#    when the content is edited/saved). We target each exact "Cas n°3.X"
#    occurrence by searching for "Cas n°" immediately followed by "3." and the
#    specific trailing digit, which is unique for each of 3.1 .. 3.4.
# ---------------------------------------------------------------------------
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-All($findText, $replaceText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null
}

Replace-All "Cas n°3.1" "Cas n°3.1"
Replace-All "Cas n°3.2" "Cas n°3.2"
Replace-All "Cas n°3.3" "Cas n°3.3"
Replace-All "Cas n°3.4" "Cas n°3.4"

# ---------------------------------------------------------------------------
# 2) "Le système fait appel au point d'extension « X »" -> "Point d'extension
#    « X »" (occurs identically for Création / Suppression / Modification).
#    The leading "Le système" run was bold; the replacement "Point" run must
#    not be bold, so after the textual replace we explicitly clear Bold on
#    the newly produced range.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.End = $d.Content.End
while ($r.Find.Execute("Le système fait appel au point", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)) {
    $r.Text = "Point"
    $r.Font.Bold = 0
    $r.Collapse(0)
}

# ---------------------------------------------------------------------------
# 3) "Le cas est appelé depuis le cas d'utilisation " -> "Au cas d'utilisation "
#    (occurs identically 3 times: in the Créer / Supprimer / Modifier use
#    cases). The third occurrence (Modifier un compte) additionally carries
#    the relocated "_GoBack" bookmark, inserted right after "Au cas".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.End = $d.Content.End
while ($r.Find.Execute("Le cas est appelé depuis le cas d’utilisation ", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)) {
    $r.Text = "Au cas d’utilisation "
    $r.Collapse(0)
}

# ---------------------------------------------------------------------------
# 4) Relocate the "_GoBack" bookmark from the end of the "Modifier un compte"
#    Ergonomie paragraph to right after "Au cas" in that same use case's
#    "Démarrage" sentence. Adding a bookmark with a name that already exists
#    moves it (Word enforces unique bookmark names), so the stale occurrence
#    at the old spot is removed automatically.
# ---------------------------------------------------------------------------
$search = $d.Content
$search.Start = 0
$search.End = $d.Content.End
$search.Find.Execute("modifier les informations d’un compte existant", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$search.Collapse(0)
$search.End = $d.Content.End
$search.Find.Execute("Au cas d’utilisation", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null

$bmPos = $search.Start + 6
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
